$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New COVID-19 case rows for 3/24/2020 (Excel date serial 43914), rows 191-275 ---

$colIndex = @{ "A" = 1; "B" = 2; "C" = 3; "D" = 4; "E" = 5; "F" = 6 }

$rows = @(
  @{ row=191; D="Wharton" },
  @{ row=192; D="Wharton" },
  @{ row=193; D="Fort Bend" },
  @{ row=194; D="Fort Bend" },
  @{ row=195; D="Fort Bend" },
  @{ row=196; D="Fort Bend" },
  @{ row=197; B="F"; C="40-50"; D="Brazoria"; E="Community Spread"; F="Y" },
  @{ row=198; B="F"; C="40-50"; D="Brazoria"; E="Community Spread" },
  @{ row=199; B="F"; C="30-40"; D="Brazoria"; E="Community Spread" },
  @{ row=200; B="F"; C="50-60"; D="Brazoria"; E="Community Spread" },
  @{ row=201; B="F"; C="30-40"; D="Brazoria"; F="Y" },
  @{ row=202; B="F"; C="30-40"; D="Brazoria"; F="Y" },
  @{ row=203; B="F"; C="20-30"; Cstyle=3; D="Brazoria"; F="Y" },
  @{ row=204; B="F"; C="70-80"; D="Brazoria" },
  @{ row=205; B="M"; C="20-30"; D="Brazoria"; F="Y" },
  @{ row=206; B="F"; C="50-60"; D="Matagorda" },
  @{ row=207; B="F"; C="80-90"; D="Matagorda" },
  @{ row=208; B="F"; C="50-60"; D="Galveston"; E="Travel" },
  @{ row=209; B="F"; C="20-30"; D="Galveston"; E="Travel" },
  @{ row=210; B="M"; C="50-60"; D="Galveston"; E="Community Spread" },
  @{ row=211; D="Houston" },
  @{ row=212; D="Houston" },
  @{ row=213; D="Houston" },
  @{ row=214; D="Houston" },
  @{ row=215; D="Houston" },
  @{ row=216; D="Houston" },
  @{ row=217; D="Houston" },
  @{ row=218; D="Houston" },
  @{ row=219; D="Houston" },
  @{ row=220; D="Houston" },
  @{ row=221; D="Houston" },
  @{ row=222; D="Houston" },
  @{ row=223; D="Houston" },
  @{ row=224; D="Houston" },
  @{ row=225; D="Houston" },
  @{ row=226; D="Houston" },
  @{ row=227; D="Houston" },
  @{ row=228; D="Houston" },
  @{ row=229; D="Houston" },
  @{ row=230; D="Houston" },
  @{ row=231; D="Houston" },
  @{ row=232; D="Houston" },
  @{ row=233; D="Houston" },
  @{ row=234; D="Houston" },
  @{ row=235; D="Houston" },
  @{ row=236; D="Houston" },
  @{ row=237; D="Houston" },
  @{ row=238; D="Houston" },
  @{ row=239; D="Houston" },
  @{ row=240; D="Houston" },
  @{ row=241; D="Houston" },
  @{ row=242; B="F"; D="Harris"; E="Exposed" },
  @{ row=243; B="F"; D="Harris"; E="Exposed" },
  @{ row=244; B="F"; D="Harris"; E="Exposed" },
  @{ row=245; B="F"; D="Harris"; E="Exposed" },
  @{ row=246; B="F"; D="Harris"; E="Travel" },
  @{ row=247; B="F"; D="Harris"; E="Travel" },
  @{ row=248; B="F"; D="Harris"; E="Community Spread" },
  @{ row=249; B="F"; D="Harris"; E="Community Spread" },
  @{ row=250; B="F"; D="Harris"; E="Community Spread" },
  @{ row=251; B="F"; D="Harris"; E="Community Spread" },
  @{ row=252; B="F"; D="Harris"; E="Community Spread" },
  @{ row=253; B="M"; D="Harris"; E="Community Spread" },
  @{ row=254; B="M"; D="Harris"; E="Community Spread" },
  @{ row=255; B="M"; D="Harris"; E="Community Spread" },
  @{ row=256; B="M"; D="Harris"; E="Community Spread" },
  @{ row=257; B="M"; D="Harris"; E="Community Spread" },
  @{ row=258; B="M"; D="Harris"; E="Community Spread" },
  @{ row=259; B="M"; D="Harris"; E="Community Spread" },
  @{ row=260; B="M"; D="Harris"; E="Community Spread" },
  @{ row=261; B="M"; D="Harris"; E="Community Spread" },
  @{ row=262; B="M"; D="Harris"; E="Community Spread" },
  @{ row=263; B="M"; D="Harris"; E="Community Spread" },
  @{ row=264; B="M"; D="Harris"; E="Community Spread" },
  @{ row=265; B="M"; D="Harris"; E="Community Spread" },
  @{ row=266; B="M"; D="Harris"; E="Community Spread" },
  @{ row=267; B="M"; C="30-40"; D="Montgomery"; E="Travel" },
  @{ row=268; B="F"; C="60-70"; D="Montgomery"; E="Community Spread" },
  @{ row=269; B="F"; C="30-40"; D="Montgomery"; E="Community Spread" },
  @{ row=270; B="M"; C="50-60"; D="Montgomery"; E="Travel" },
  @{ row=271; D="Brazos" },
  @{ row=272; D="Brazos" },
  @{ row=273; D="Brazos" },
  @{ row=274; D="Brazos" },
  @{ row=275; B="F"; C="50-60"; D="Brazoria" }
)

# Column A on every new row is the report date (3/24/2020), formatted like the
# existing date column (copy the format from the last existing date cell so we
# reuse the existing style instead of minting new cellXfs entries).
$firstRow = $rows[0]["row"]
$lastRow = $rows[$rows.Count - 1]["row"]
$ws.Range("A190").Copy($ws.Range("A" + $firstRow + ":A" + $lastRow))
$ws.Range("A" + $firstRow + ":A" + $lastRow).Value = 43914

foreach ($rowData in $rows) {
  $r = $rowData["row"]
  foreach ($col in @("B", "C", "D", "E", "F")) {
    if ($rowData.ContainsKey($col)) {
      $ci = $colIndex[$col]
      $cell = $ws.Cells.Item($r, $ci)
      $cell.Value = $rowData[$col]

      # One special case in the source data: C203 ("20-30") carries a stray
      # date-ish number format (d-mmm / built-in numFmtId 16) instead of the
      # default General format - replicate that exact cell-level formatting.
      $styleKey = $col + "style"
      if ($rowData.ContainsKey($styleKey)) {
        $cell.NumberFormat = "d-mmm"
      }
    }
  }
}

# Leave the selection where the author left it when they saved.
$ws.Range("D253").Select()
